$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before "总计" (so it becomes the 6th
#    tab, pushing "总计" to the 7th / last position).
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1Sheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$q1Sheet.Name = "2022-Q1"

# Header row (matches the style used on the other quarterly sheets: bold,
# centered, thin border).
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

$q1Header = $q1Sheet.Range("B1:H1")
$q1Header.Font.Bold = $true
$q1Header.HorizontalAlignment = -4108
$q1Header.VerticalAlignment = -4160
$q1Header.Borders.LineStyle = 1

# Row-index cell in column A uses the same bold/centered/bordered style.
$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("A2").Font.Bold = $true
$q1Sheet.Range("A2").HorizontalAlignment = -4108
$q1Sheet.Range("A2").VerticalAlignment = -4160
$q1Sheet.Range("A2").Borders.LineStyle = 1

# Data row - the numeric-looking values (fund code, scale, position, etc.)
# are stored as text on the sibling sheets, so force text with a leading
# apostrophe just like the existing quarterly sheets do.
$q1Sheet.Range("B2").Value = "'161225"
$q1Sheet.Range("C2").Value = "国投瑞银瑞盈灵活配置混合（LOF）"
$q1Sheet.Range("D2").Value = "'2.28"
$q1Sheet.Range("E2").Value = "'94.46"
$q1Sheet.Range("F2").Value = "'3.67"
$q1Sheet.Range("G2").Value = "'0.0837"
$q1Sheet.Range("H2").Value = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row right under the
#    header for the "2022-Q1" summary, shifting the existing rows down and
#    re-numbering the column-A row index.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A2").Font.Bold = $true
$totalSheet.Range("A2").HorizontalAlignment = -4108
$totalSheet.Range("A2").VerticalAlignment = -4160
$totalSheet.Range("A2").Borders.LineStyle = 1

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.08

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
